$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Added standard quantity for 4 servings - update ingredient quantities
$ws.Range("C2").Value = 400
$ws.Range("C3").Value = 100
$ws.Range("C4").Value = 500
$ws.Range("C5").Value = 300
$ws.Range("C6").Value = 6
$ws.Range("C7").Value = 6
$ws.Range("C8").Value = 6
$ws.Range("C9").Value = 4
$ws.Range("C10").Value = 8
$ws.Range("C11").Value = 400
$ws.Range("C12").Value = 1
$ws.Range("C13").Value = 200
$ws.Range("C14").Value = 2
$ws.Range("C15").Value = 1
$ws.Range("C16").Value = 80

# Fix formatting on the last row's Unit cell to match the rest of the column
# (right aligned, like D2:D15)
$ws.Range("D16").HorizontalAlignment = -4152

# Update the saved cursor/selection position
$ws.Range("C8").Select()
